$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Stata"

# --- Headers (row 1) ---
$ws.Range("D1").Value = "Laptop"
$ws.Range("E1").Value = "Currency"
$ws.Range("F1").Value = "Exchange rate"
$ws.Range("G1").Value = "Laptop USD"
$ws.Range("H1").Value = "Percent license"
$ws.Range("I1").Value = "Percent laptop"

# --- Move the old "percent" formula from column D to column H ---
$ws.Range("H2").Formula = "=C2/B2*100"
$ws.Range("H3").Formula = "=C3/B3*100"
$ws.Range("H4").Formula = "=C4/B4*100"

# --- New data: Laptop local price, currency, exchange rate ---
$ws.Range("D2").Value = 849
$ws.Range("E2").Value = "USD"
$ws.Range("F2").Value = 1

$ws.Range("D3").Value = 799
$ws.Range("E3").Value = "EUR"
$ws.Range("F3").Value = 1.17

$ws.Range("D4").Value = 18440000
$ws.Range("E4").Value = "VND"
$ws.Range("F4").Value = 38.32

# --- Laptop USD conversion formulas ---
$ws.Range("G2").Formula = "=D2"
$ws.Range("G3").Formula = "=D3*F3"
$ws.Range("G4").Formula = "=D4/1000000*F4"

# --- Percent laptop formulas ---
$ws.Range("I2").Formula = "=G2/B2*100"
$ws.Range("I3").Formula = "=G3/B3*100"
$ws.Range("I4").Formula = "=G4/B4*100"

# --- Number formats ---
$ws.Range("G2:G4").NumberFormat = "[$$-409]#,##0.00;[RED]\-[$$-409]#,##0.00"
$ws.Range("H2:H4").NumberFormat = "General"
$ws.Range("I2:I4").NumberFormat = "General"

# --- Column widths (closest achievable to the target character widths,
# Excel's ColumnWidth is quantized to a pixel grid before being stored) ---
$ws.Columns.Item(4).ColumnWidth = 13.43
$ws.Columns.Item(7).ColumnWidth = 9.92

# --- Selection to match target ---
$ws.Range("I4").Select() | Out-Null
